$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cv_entries")

# --- Row 2: Ph.D. entry - GPA updated from "4+" to "4.12" ---
$ws.Range("G2").Value = "GPA: 4.12`nSpecialization in Educational Data Science with R`nAdvisor: Prof. David Liebowitz"

# --- Row 3: M.S. Economics entry - GPA updated from "4+" to "4.12" ---
$ws.Range("G3").Value = "GPA: 4.12`nAdvisor: Prof. Glen Waddell`nThesis Title: Student First, or Budget First? Examining School Districts" + [char]0x2019 + " Spending Trends for Special Education Students."

# --- Row 5: Research Assistant (David Liebowitz) - description reworded ---
$ws.Range("G5").Value = "Worked with Dr. David Liebowitz on research projects exploring principle effects and school discipline"

# --- Row 6: Research Assistant (Gina Biancarosa / Patrick Kennedy) - description reworded ---
$ws.Range("G6").Value = "Worked with Dr. Gina Biancarosa and Dr. Patrick Kennedy on research projects related to Dynamic Indicators of Basic Early Literacy Skill (DIBELS8) reading assessments"

# --- Row 4: B.S. (H) Mathematics entry - additional_info filled in ---
$ws.Range("G4").Value = "Composite Percentage: 87% (equivalent to a GPA of 4)"

# --- Update the sheet view: scroll position, zoom, and selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("C1").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G4").Select()
